$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entries for 2014-11-29 and 2014-11-30 (rows 20 and 21),
# inserted right after the existing last entry (row 19, 2014-11-26).

# Row 20: 2014-11-29, 13:00 -> 19:00
$ws.Cells.Item(20, 1).Value = 41972
$ws.Cells.Item(20, 2).Value = 0.541666666666667
$ws.Cells.Item(20, 3).Value = 0.791666666666667
$ws.Cells.Item(20, 4).Formula = "=ROUND(ABS(C20-B20) * 24, 1)"

# Row 21: 2014-11-30, 01:00 -> 05:30
$ws.Cells.Item(21, 1).Value = 41973
$ws.Cells.Item(21, 2).Value = 0.0416666666666667
$ws.Cells.Item(21, 3).Value = 0.229166666666667
$ws.Cells.Item(21, 4).Formula = "=ROUND(ABS(C21-B21) * 24, 1)"

# Match the formatting (date / time number formats) used by the rest of
# the table by copying it down from the row directly above the new ones.
$ws.Range("A19:D19").Copy()
$ws.Range("A20:D21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to D20, mirroring where the author was
# working after adding the new rows.
[void]$ws.Range("D20").Select()
